{"js": "// SMARTNODES.docx edit:\n//  1) Renumber the \"smarthosting\" bookmark from id 1 -> 0 (Word does this\n//     automatically when bookmarks are (re)created; deleting and\n//     re-inserting the sole bookmark in the document yields id 0).\n//  2) In the \"100,000 Smart:\" bullet, drop the two exchange hyperlinks\n//     (CryptoBridge, HitBTC) -- including their HYPERLINK field codes --\n//     and the now-dangling \" such as \" / \", \" connective text, leaving\n//     \"...obtained from exchanges.\" followed by the untouched rest of the\n//     sentence (the \"here\" hyperlink, etc.).\n\nconst doc = context.document;\n\n// --- Change 1: bookmark id 1 -> 0 -------------------------------------\nconst bookmarkRange = doc.getBookmarkRange(\"smarthosting\");\ndoc.deleteBookmark(\"smarthosting\");\nbookmarkRange.insertBookmark(\"smarthosting\");\nawait context.sync();\n\n// --- Change 2: remove the CryptoBridge / HitBTC exchange mentions ------\nconst paragraphs = doc.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nparagraphs.items.forEach((p) => p.load(\"text\"));\nawait context.sync();\n\nconst target = paragraphs.items.find(\n  (p) => p.text.indexOf(\"100,000 Smart:\") !== -1\n);\n\nif (target) {\n  const paraRange = target.getRange();\n  const fields = paraRange.fields;\n  fields.load(\"items\");\n  await context.sync();\n\n  fields.items.forEach((f) => f.load(\"code\"));\n  await context.sync();\n\n  // Delete only the HYPERLINK fields pointing at the two exchanges; leave\n  // the \"here\" hyperlink field (and everything else) untouched.\n  fields.items\n    .filter((f) => {\n      const code = f.code || \"\";\n      return code.indexOf(\"crypto-bridge.org\") !== -1 || code.indexOf(\"hitbtc.com\") !== -1;\n    })\n    .forEach((f) => f.delete());\n  await context.sync();\n\n  // With the two fields gone, the remaining dangling connective text reads\n  // \"... exchanges such as , . For the full list ...\" -- strip it back to\n  // \"... exchanges. For the full list ...\".\n  const dangling = target.search(\" such as , \", { matchCase: true });\n  dangling.load(\"items\");\n  await context.sync();\n\n  if (dangling.items.length > 0) {\n    dangling.items[0].delete();\n    await context.sync();\n  }\n}\n", "ps1": "# SMARTNODES.docx edit:\n#  1) Renumber the \"smarthosting\" bookmark from id 1 -> 0 (Word assigns\n#     this automatically; deleting and re-adding the sole bookmark at the\n#     same Range yields id 0).\n#  2) In the \"100,000 Smart:\" bullet, drop the two exchange hyperlinks\n#     (CryptoBridge, HitBTC) -- including their HYPERLINK field codes --\n#     and the now-dangling \" such as \" / \", \" connective text, leaving\n#     \"...obtained from exchanges.\" followed by the untouched rest of the\n#     sentence (the \"here\" hyperlink, etc.).\n\n$d = $word.ActiveDocument\n\n# --- Change 1: bookmark id 1 -> 0 --------------------------------------\n$bm = $d.Bookmarks(\"smarthosting\")\n$bmRange = $bm.Range\n$d.Bookmarks(\"smarthosting\").Delete()\n$d.Bookmarks.Add(\"smarthosting\", $bmRange)\n\n# --- Change 2: remove the CryptoBridge / HitBTC exchange mentions ------\n$fieldsToDelete = @()\nforeach ($f in $d.Fields) {\n    $code = $f.Code.Text\n    if ($code -like \"*crypto-bridge.org*\" -or $code -like \"*hitbtc.com*\") {\n        $fieldsToDelete += $f\n    }\n}\nforeach ($f in $fieldsToDelete) {\n    $f.Delete()\n}\n\n# With the two fields gone, the remaining dangling connective text reads\n# \"... exchanges such as , . For the full list ...\" -- strip it back to\n# \"... exchanges. For the full list ...\".\n$findRange = $d.Content\n$findRange.Find.ClearFormatting()\n$found = $findRange.Find.Execute(\" such as , \")\nif ($found) {\n    $findRange.Text = \"\"\n}\n"}
